$d = $word.ActiveDocument

# --- 1) "The target file name will have the target language as a suffix."
#        -> "...language code as a suffix." split across 3 runs:
#           "...language " | "code" | " as a suffix."
$r = $d.Content
$found = $r.Find.Execute("target language ")
$r.Collapse(0)
$start = $r.Start

# Insert the new word "code" right after "target language ".
$insCode = $d.Range($start, $start)
$insCode.InsertAfter("code")

# Insert the separating space before "as a suffix." (currently directly
# abutting "code" after the insertion above).
$insSpace = $d.Range($start + 4, $start + 4)
$insSpace.InsertAfter(" ")

# Re-apply (no-op) character formatting on just the "code" span so it is
# materialized as its own run, splitting it off from the surrounding text
# that keeps the original formatting.
$codeRange = $d.Range($start, $start + 4)
$codeRange.Font.Name = "Cascadia Mono"
$codeRange.Font.Size = 11

# --- 2) Give the two trailing empty runs (previously <w:rPr/>) the same
#        Cascadia Mono / 11pt formatting as the rest of the paragraphs.
$p1 = $d.Paragraphs.Item(13)
$p1.Range.Font.Name = "Cascadia Mono"
$p1.Range.Font.Size = 11
$p1.Range.Font.SizeBi = 11

$p2 = $d.Paragraphs.Item(15)
$p2.Range.Font.Name = "Cascadia Mono"
$p2.Range.Font.Size = 11
$p2.Range.Font.SizeBi = 11
